# edit.ps1
# Commit: feat: add 2022-Q1 data
#
# Inserts a new "2022-Q1" fund-holdings detail sheet (same layout as the other
# quarterly tabs: 基金代码/基金名称/基金规模/股票总仓位/仓位占比/持有市值(亿元)/仓位排名)
# between the existing "2021-Q4" tab and the "总计" summary tab, and refreshes the
# "总计" sheet with a new leading row summarising the 2022-Q1 data.

$wb = $excel.ActiveWorkbook
$detailTemplate = $wb.Worksheets.Item("2021-Q4")
$total = $wb.Worksheets.Item("总计")

# ------------------------------------------------------------------
# 1. New detail sheet "2022-Q1", inserted right before "总计"
# ------------------------------------------------------------------
$q1 = $wb.Worksheets.Add($total)
$q1.Name = "2022-Q1"

# Borrow the bold/bordered header + index-column formatting from the 2021-Q4 sheet
# (B1:H1 header row, A-column row index cells) so the new tab matches its siblings.
$detailTemplate.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$detailTemplate.Range("A2").Copy()
$q1.Range("A2:A30").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header row
$q1.Range('B1').Value = '基金代码'
$q1.Range('C1').Value = '基金名称'
$q1.Range('D1').Value = '基金规模'
$q1.Range('E1').Value = '股票总仓位'
$q1.Range('F1').Value = '仓位占比'
$q1.Range('G1').Value = '持有市值(亿元)'
$q1.Range('H1').Value = '仓位排名'

# Columns B-G hold text values in the source data (fund codes need their leading
# zeros, and the numeric-looking figures are plain text) - force text format first
# so Excel doesn't silently convert them to numbers.
$q1.Range("B2:G30").NumberFormat = "@"

# Data rows (A=index, B=code, C=name, D=size, E=position, F=share, G=value, H=rank)
$q1.Cells.Item(2, 1).Value = 0
$q1.Cells.Item(2, 2).Value = '002121'
$q1.Cells.Item(2, 3).Value = '广发沪港深新起点股票A'
$q1.Cells.Item(2, 4).Value = '34.56'
$q1.Cells.Item(2, 5).Value = '91.46'
$q1.Cells.Item(2, 6).Value = '6.10'
$q1.Cells.Item(2, 7).Value = '2.1082'
$q1.Cells.Item(2, 8).Value = 6
$q1.Cells.Item(3, 1).Value = 1
$q1.Cells.Item(3, 2).Value = '010041'
$q1.Cells.Item(3, 3).Value = '嘉实港股优势混合A'
$q1.Cells.Item(3, 4).Value = '49.98'
$q1.Cells.Item(3, 5).Value = '92.41'
$q1.Cells.Item(3, 6).Value = '3.69'
$q1.Cells.Item(3, 7).Value = '1.8443'
$q1.Cells.Item(3, 8).Value = 8
$q1.Cells.Item(4, 1).Value = 2
$q1.Cells.Item(4, 2).Value = '001878'
$q1.Cells.Item(4, 3).Value = '嘉实沪港深精选股票'
$q1.Cells.Item(4, 4).Value = '23.17'
$q1.Cells.Item(4, 5).Value = '93.29'
$q1.Cells.Item(4, 6).Value = '3.86'
$q1.Cells.Item(4, 7).Value = '0.8944'
$q1.Cells.Item(4, 8).Value = 7
$q1.Cells.Item(5, 1).Value = 3
$q1.Cells.Item(5, 2).Value = '010761'
$q1.Cells.Item(5, 3).Value = '华商甄选回报混合'
$q1.Cells.Item(5, 4).Value = '20.63'
$q1.Cells.Item(5, 5).Value = '93.93'
$q1.Cells.Item(5, 6).Value = '3.93'
$q1.Cells.Item(5, 7).Value = '0.8108'
$q1.Cells.Item(5, 8).Value = 6
$q1.Cells.Item(6, 1).Value = 4
$q1.Cells.Item(6, 2).Value = '011856'
$q1.Cells.Item(6, 3).Value = '安信均衡成长18个月持有期混合型证券投资基金A'
$q1.Cells.Item(6, 4).Value = '6.52'
$q1.Cells.Item(6, 5).Value = '89.66'
$q1.Cells.Item(6, 6).Value = '10.31'
$q1.Cells.Item(6, 7).Value = '0.6722'
$q1.Cells.Item(6, 8).Value = 1
$q1.Cells.Item(7, 1).Value = 5
$q1.Cells.Item(7, 2).Value = '009715'
$q1.Cells.Item(7, 3).Value = '汇添富策略增长两年封闭运作灵活配置混合'
$q1.Cells.Item(7, 4).Value = '11.81'
$q1.Cells.Item(7, 5).Value = '88.41'
$q1.Cells.Item(7, 6).Value = '4.34'
$q1.Cells.Item(7, 7).Value = '0.5126'
$q1.Cells.Item(7, 8).Value = 4
$q1.Cells.Item(8, 1).Value = 6
$q1.Cells.Item(8, 2).Value = '003304'
$q1.Cells.Item(8, 3).Value = '前海开源沪港深核心资源灵活配置混合A'
$q1.Cells.Item(8, 4).Value = '5.91'
$q1.Cells.Item(8, 5).Value = '93.10'
$q1.Cells.Item(8, 6).Value = '7.37'
$q1.Cells.Item(8, 7).Value = '0.4356'
$q1.Cells.Item(8, 8).Value = 8
$q1.Cells.Item(9, 1).Value = 7
$q1.Cells.Item(9, 2).Value = '008891'
$q1.Cells.Item(9, 3).Value = '安信价值成长混合A'
$q1.Cells.Item(9, 4).Value = '2.83'
$q1.Cells.Item(9, 5).Value = '92.41'
$q1.Cells.Item(9, 6).Value = '10.18'
$q1.Cells.Item(9, 7).Value = '0.2881'
$q1.Cells.Item(9, 8).Value = 1
$q1.Cells.Item(10, 1).Value = 8
$q1.Cells.Item(10, 2).Value = '009880'
$q1.Cells.Item(10, 3).Value = '安信成长动力一年持有期混合'
$q1.Cells.Item(10, 4).Value = '2.75'
$q1.Cells.Item(10, 5).Value = '93.58'
$q1.Cells.Item(10, 6).Value = '9.49'
$q1.Cells.Item(10, 7).Value = '0.2610'
$q1.Cells.Item(10, 8).Value = 2
$q1.Cells.Item(11, 1).Value = 9
$q1.Cells.Item(11, 2).Value = '010042'
$q1.Cells.Item(11, 3).Value = '嘉实港股优势混合C'
$q1.Cells.Item(11, 4).Value = '5.78'
$q1.Cells.Item(11, 5).Value = '92.41'
$q1.Cells.Item(11, 6).Value = '3.69'
$q1.Cells.Item(11, 7).Value = '0.2133'
$q1.Cells.Item(11, 8).Value = 8
$q1.Cells.Item(12, 1).Value = 10
$q1.Cells.Item(12, 2).Value = '003305'
$q1.Cells.Item(12, 3).Value = '前海开源沪港深核心资源灵活配置混合C'
$q1.Cells.Item(12, 4).Value = '2.19'
$q1.Cells.Item(12, 5).Value = '93.10'
$q1.Cells.Item(12, 6).Value = '7.37'
$q1.Cells.Item(12, 7).Value = '0.1614'
$q1.Cells.Item(12, 8).Value = 8
$q1.Cells.Item(13, 1).Value = 11
$q1.Cells.Item(13, 2).Value = '008488'
$q1.Cells.Item(13, 3).Value = '华商恒益稳健混合'
$q1.Cells.Item(13, 4).Value = '2.03'
$q1.Cells.Item(13, 5).Value = '58.93'
$q1.Cells.Item(13, 6).Value = '5.79'
$q1.Cells.Item(13, 7).Value = '0.1175'
$q1.Cells.Item(13, 8).Value = 1
$q1.Cells.Item(14, 1).Value = 12
$q1.Cells.Item(14, 2).Value = '014746'
$q1.Cells.Item(14, 3).Value = '贝莱德港股通远景视野混合A'
$q1.Cells.Item(14, 4).Value = '5.05'
$q1.Cells.Item(14, 5).Value = '53.79'
$q1.Cells.Item(14, 6).Value = '2.26'
$q1.Cells.Item(14, 7).Value = '0.1141'
$q1.Cells.Item(14, 8).Value = 6
$q1.Cells.Item(15, 1).Value = 13
$q1.Cells.Item(15, 2).Value = '008892'
$q1.Cells.Item(15, 3).Value = '安信价值成长混合C'
$q1.Cells.Item(15, 4).Value = '0.62'
$q1.Cells.Item(15, 5).Value = '92.41'
$q1.Cells.Item(15, 6).Value = '10.18'
$q1.Cells.Item(15, 7).Value = '0.0631'
$q1.Cells.Item(15, 8).Value = 1
$q1.Cells.Item(16, 1).Value = 14
$q1.Cells.Item(16, 2).Value = '241001'
$q1.Cells.Item(16, 3).Value = '华宝海外中国混合(QDII)'
$q1.Cells.Item(16, 4).Value = '0.83'
$q1.Cells.Item(16, 5).Value = '86.89'
$q1.Cells.Item(16, 6).Value = '6.28'
$q1.Cells.Item(16, 7).Value = '0.0521'
$q1.Cells.Item(16, 8).Value = 4
$q1.Cells.Item(17, 1).Value = 15
$q1.Cells.Item(17, 2).Value = '014747'
$q1.Cells.Item(17, 3).Value = '贝莱德港股通远景视野混合C'
$q1.Cells.Item(17, 4).Value = '2.23'
$q1.Cells.Item(17, 5).Value = '53.79'
$q1.Cells.Item(17, 6).Value = '2.26'
$q1.Cells.Item(17, 7).Value = '0.0504'
$q1.Cells.Item(17, 8).Value = 6
$q1.Cells.Item(18, 1).Value = 16
$q1.Cells.Item(18, 2).Value = '008253'
$q1.Cells.Item(18, 3).Value = '华宝致远混合（QDII）A'
$q1.Cells.Item(18, 4).Value = '0.70'
$q1.Cells.Item(18, 5).Value = '85.00'
$q1.Cells.Item(18, 6).Value = '6.31'
$q1.Cells.Item(18, 7).Value = '0.0442'
$q1.Cells.Item(18, 8).Value = 2
$q1.Cells.Item(19, 1).Value = 17
$q1.Cells.Item(19, 2).Value = '010024'
$q1.Cells.Item(19, 3).Value = '广发沪港深新起点股票C'
$q1.Cells.Item(19, 4).Value = '0.72'
$q1.Cells.Item(19, 5).Value = '91.46'
$q1.Cells.Item(19, 6).Value = '6.10'
$q1.Cells.Item(19, 7).Value = '0.0439'
$q1.Cells.Item(19, 8).Value = 6
$q1.Cells.Item(20, 1).Value = 18
$q1.Cells.Item(20, 2).Value = '009017'
$q1.Cells.Item(20, 3).Value = '银华港股通精选股票'
$q1.Cells.Item(20, 4).Value = '0.91'
$q1.Cells.Item(20, 5).Value = '86.12'
$q1.Cells.Item(20, 6).Value = '4.50'
$q1.Cells.Item(20, 7).Value = '0.0410'
$q1.Cells.Item(20, 8).Value = 9
$q1.Cells.Item(21, 1).Value = 19
$q1.Cells.Item(21, 2).Value = '012924'
$q1.Cells.Item(21, 3).Value = '华夏新时代灵活配置混合（QDII）美元现汇'
$q1.Cells.Item(21, 4).Value = '2.56'
$q1.Cells.Item(21, 5).Value = '84.71'
$q1.Cells.Item(21, 6).Value = '1.57'
$q1.Cells.Item(21, 7).Value = '0.0402'
$q1.Cells.Item(21, 8).Value = 10
$q1.Cells.Item(22, 1).Value = 20
$q1.Cells.Item(22, 2).Value = '012925'
$q1.Cells.Item(22, 3).Value = '华夏新时代灵活配置混合（QDII）美元现钞'
$q1.Cells.Item(22, 4).Value = '2.56'
$q1.Cells.Item(22, 5).Value = '84.71'
$q1.Cells.Item(22, 6).Value = '1.57'
$q1.Cells.Item(22, 7).Value = '0.0402'
$q1.Cells.Item(22, 8).Value = 10
$q1.Cells.Item(23, 1).Value = 21
$q1.Cells.Item(23, 2).Value = '011857'
$q1.Cells.Item(23, 3).Value = '安信均衡成长18个月持有期混合型证券投资基金C'
$q1.Cells.Item(23, 4).Value = '0.32'
$q1.Cells.Item(23, 5).Value = '89.66'
$q1.Cells.Item(23, 6).Value = '10.31'
$q1.Cells.Item(23, 7).Value = '0.0330'
$q1.Cells.Item(23, 8).Value = 1
$q1.Cells.Item(24, 1).Value = 22
$q1.Cells.Item(24, 2).Value = '005701'
$q1.Cells.Item(24, 3).Value = '上投摩根香港精选港股通混合'
$q1.Cells.Item(24, 4).Value = '0.44'
$q1.Cells.Item(24, 5).Value = '84.37'
$q1.Cells.Item(24, 6).Value = '3.50'
$q1.Cells.Item(24, 7).Value = '0.0154'
$q1.Cells.Item(24, 8).Value = 4
$q1.Cells.Item(25, 1).Value = 23
$q1.Cells.Item(25, 2).Value = '501303'
$q1.Cells.Item(25, 3).Value = '广发港股通恒生综合中型股指数(LOF)A'
$q1.Cells.Item(25, 4).Value = '0.34'
$q1.Cells.Item(25, 5).Value = '92.39'
$q1.Cells.Item(25, 6).Value = '2.29'
$q1.Cells.Item(25, 7).Value = '0.0078'
$q1.Cells.Item(25, 8).Value = 2
$q1.Cells.Item(26, 1).Value = 24
$q1.Cells.Item(26, 2).Value = '008254'
$q1.Cells.Item(26, 3).Value = '华宝致远混合（QDII）C'
$q1.Cells.Item(26, 4).Value = '0.12'
$q1.Cells.Item(26, 5).Value = '85.00'
$q1.Cells.Item(26, 6).Value = '6.31'
$q1.Cells.Item(26, 7).Value = '0.0076'
$q1.Cells.Item(26, 8).Value = 2
$q1.Cells.Item(27, 1).Value = 25
$q1.Cells.Item(27, 2).Value = '501067'
$q1.Cells.Item(27, 3).Value = '招商富时中国A-H50指数（LOF）A'
$q1.Cells.Item(27, 4).Value = '0.21'
$q1.Cells.Item(27, 5).Value = '94.63'
$q1.Cells.Item(27, 6).Value = '2.85'
$q1.Cells.Item(27, 7).Value = '0.0060'
$q1.Cells.Item(27, 8).Value = 9
$q1.Cells.Item(28, 1).Value = 26
$q1.Cells.Item(28, 2).Value = '004996'
$q1.Cells.Item(28, 3).Value = '广发港股通恒生综合中型股指数(LOF)C'
$q1.Cells.Item(28, 4).Value = '0.11'
$q1.Cells.Item(28, 5).Value = '92.39'
$q1.Cells.Item(28, 6).Value = '2.29'
$q1.Cells.Item(28, 7).Value = '0.0025'
$q1.Cells.Item(28, 8).Value = 2
$q1.Cells.Item(29, 1).Value = 27
$q1.Cells.Item(29, 2).Value = '160922'
$q1.Cells.Item(29, 3).Value = '大成恒生综合中小型股指数(QDII-LOF)A'
$q1.Cells.Item(29, 4).Value = '0.10'
$q1.Cells.Item(29, 5).Value = '92.44'
$q1.Cells.Item(29, 6).Value = '1.78'
$q1.Cells.Item(29, 7).Value = '0.0018'
$q1.Cells.Item(29, 8).Value = 2
$q1.Cells.Item(30, 1).Value = 28
$q1.Cells.Item(30, 2).Value = '501068'
$q1.Cells.Item(30, 3).Value = '招商富时中国A-H50指数（LOF）C'
$q1.Cells.Item(30, 4).Value = '0.05'
$q1.Cells.Item(30, 5).Value = '94.63'
$q1.Cells.Item(30, 6).Value = '2.85'
$q1.Cells.Item(30, 7).Value = '0.0014'
$q1.Cells.Item(30, 8).Value = 9

Write-Host ("2022-Q1 sheet populated: " + $q1.Range("B2").Value + " / " + $q1.Cells.Item(30, 2).Value)

# ------------------------------------------------------------------
# 2. Refresh "总计": shift existing rows down one and add the new
#    2022-Q1 summary row at the top (row 2)
# ------------------------------------------------------------------
# Row 7 is brand new (sheet used to stop at row 6) - clone the bold/bordered
# index-column style from A2 before writing into it.
$total.Range("A2").Copy()
$total.Range("A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = '2022-Q1'
$total.Cells.Item(2, 3).Value = 29
$total.Cells.Item(2, 4).Value = 8.880000000000001
$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(3, 2).Value = '2021-Q4'
$total.Cells.Item(3, 3).Value = 16
$total.Cells.Item(3, 4).Value = 4.43
$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(4, 2).Value = '2021-Q3'
$total.Cells.Item(4, 3).Value = 19
$total.Cells.Item(4, 4).Value = 7.49
$total.Cells.Item(5, 1).Value = 3
$total.Cells.Item(5, 2).Value = '2021-Q2'
$total.Cells.Item(5, 3).Value = 16
$total.Cells.Item(5, 4).Value = 10.96
$total.Cells.Item(6, 1).Value = 4
$total.Cells.Item(6, 2).Value = '2021-Q1'
$total.Cells.Item(6, 3).Value = 22
$total.Cells.Item(6, 4).Value = 13.83
$total.Cells.Item(7, 1).Value = 5
$total.Cells.Item(7, 2).Value = '2020-Q4'
$total.Cells.Item(7, 3).Value = 20
$total.Cells.Item(7, 4).Value = 10.92

Write-Host ("总计 refreshed: " + $total.Range("B2").Value + " now leads the summary")

Write-Host "Final sheet order:"
foreach ($s in $wb.Worksheets) { Write-Host (" - " + $s.Name) }
